$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")
$model = $wb.Worksheets.Item("Model")

# --- Sheet view: scroll/zoom on Main ---
$ws.Application.ActiveWindow.Zoom = 160
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Application.ActiveWindow.ScrollColumn = 1

# --- L2 price update ---
$ws.Range("L2").Value = 219.94

# --- Row 3: clear stray borders/fills from label cells, drop F3 entirely ---
$ws.Range("C3").ClearFormats()
$ws.Range("D3").ClearFormats()
$ws.Range("E3").ClearFormats()
$ws.Range("F3").Clear()
$ws.Range("G3").ClearFormats()

# L3 becomes a formula
$ws.Range("L3").Formula = "=1349.64018/13"
# M3: change referenced quarter label from Q121 to Q322
$ws.Range("M3").Value = "Q322"

# --- Row 4 ---
$ws.Range("C4").ClearFormats()
$ws.Range("D4").ClearFormats()
$ws.Range("E4").ClearFormats()
$ws.Range("F4").ClearFormats()
$ws.Range("G4").ClearFormats()
# L4 keeps its formula (recalculates automatically since L2/L3 changed)

# --- Row 5 ---
$ws.Range("C5").ClearFormats()
$ws.Range("D5").Clear()
$ws.Range("E5").Clear()
$ws.Range("F5").Clear()
$ws.Range("G5").ClearFormats()
$ws.Range("L5").Formula = "=4197.132+871.998"
$ws.Range("M5").Value = "Q322"

# --- Row 6 ---
$ws.Range("C6").Clear()
$ws.Range("D6").Clear()
$ws.Range("E6").Clear()
$ws.Range("F6").Clear()
$ws.Range("G6").ClearFormats()
$ws.Range("L6").Formula = "=441.275+208"
$ws.Range("M6").Value = "Q322"

# --- Row 7 ---
$ws.Range("C7").Clear()
$ws.Range("D7").Clear()
$ws.Range("E7").Clear()
$ws.Range("F7").Clear()
$ws.Range("G7").Clear()
# L7 keeps its formula (=L4-L5+L6), recalculates automatically

# --- Row 10 ---
$ws.Range("C10").ClearFormats()
$ws.Range("D10").Clear()
$ws.Range("E10").ClearFormats()
$ws.Range("F10").Clear()
$ws.Range("G10").ClearFormats()

# --- Row 11 ---
$ws.Range("C11").Clear()
$ws.Range("D11").Clear()
$ws.Range("E11").ClearFormats()
$ws.Range("F11").Clear()
$ws.Range("G11").Clear()

# --- Row 12 ---
$ws.Range("C12").ClearFormats()
$ws.Range("D12").Clear()
$ws.Range("E12").Clear()
$ws.Range("F12").Clear()
$ws.Range("G12").Clear()

# --- Row 13 ---
$ws.Range("C13").Clear()
$ws.Range("D13").Clear()
$ws.Range("E13").ClearFormats()
$ws.Range("F13").Clear()
$ws.Range("G13").Clear()

# --- Row 14 ---
$ws.Range("C14").Clear()
$ws.Range("D14").Clear()
$ws.Range("E14").ClearFormats()
$ws.Range("F14").Clear()
$ws.Range("G14").Clear()

# --- Row 15 ---
$ws.Range("C15").Clear()
$ws.Range("D15").Clear()
$ws.Range("E15").ClearFormats()
$ws.Range("F15").Clear()
$ws.Range("G15").Clear()

# --- Row 16 ---
$ws.Range("C16").Clear()
$ws.Range("D16").Clear()
$ws.Range("E16").ClearFormats()
$ws.Range("F16").Clear()
$ws.Range("G16").Clear()

# --- Row 17 ---
$ws.Range("C17").Clear()
$ws.Range("D17").Clear()
$ws.Range("E17").ClearFormats()
$ws.Range("F17").Clear()
$ws.Range("G17").Clear()

# --- Row 18 ---
$ws.Range("C18").ClearFormats()
$ws.Range("D18").Clear()
$ws.Range("E18").ClearFormats()
$ws.Range("F18").Clear()
$ws.Range("G18").Clear()
